$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (target stored width = 23)
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668

# Insert a new row at position 5 - shifts the old summary block (rows 5-14) down to
# rows 6-15, and creates a new blank row 5 that we'll fill with the 4th line item.
$ws.Range("A5").EntireRow.Insert()

# --- Line items (rows 2-5) ---
$ws.Range("A2").Value = "Licencia Excel"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 22380

$ws.Range("A3").Value = "Lista Negra"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1200
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 17904

$ws.Range("A4").Value = "Ingeniero Informático"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1500
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 12000

$ws.Range("A5").Value = "Licencia Bizagui"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1119

# --- Summary block (rows 7-15, row 6 stays blank) ---
$ws.Range("D7").Value = "TOTAL"
$ws.Range("E7").Value = 53403

$ws.Range("D8").Value = "Reserva de contingencia"
$ws.Range("E8").Value = 1000

$ws.Range("D9").Value = "Linea Base de Costos"
$ws.Range("E9").Value = 54403
$ws.Range("F9").ClearContents()

$ws.Range("D10").Value = "Reserva de gestion"
$ws.Range("E10").Value = "PV"
$ws.Range("F10").Value = 0.05

$ws.Range("D11").Value = "Presupuesto"
$ws.Range("E11").Value = 50000

$ws.Range("D12").Value = "Ganancia"
$ws.Range("E12").Value = "PV"
$ws.Range("F12").Value = 0.1

$ws.Range("D13").Value = "Total con ganancia"
$ws.Range("E13").Value = 104403

$ws.Range("D14").Value = "IGV"
$ws.Range("E14").Value = "PV"
$ws.Range("F14").Value = 0.04

$ws.Range("D15").Value = "Total"
$ws.Range("E15").Value = 104403.04
